$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 26 results: row 36 = MI vs RR, row 37 = DC vs KKR, row 38 = PBKS vs RCB
# Player point columns (in player order): E, H, K, N, Q, T, W, Z, AC

$data = @{
    36 = @{ "E" = 80;  "H" = 70;  "K" = 60; "N" = 30; "Q" = 50; "T" = 0;  "W" = 100; "Z" = 40; "AC" = 20  }
    37 = @{ "E" = 80;  "H" = 50;  "K" = 0;  "N" = 40; "Q" = 30; "T" = 60; "W" = 100; "Z" = 20; "AC" = 70  }
    38 = @{ "E" = 80;  "H" = 70;  "K" = 30; "N" = 60; "Q" = 50; "T" = 20; "W" = 40;  "Z" = 0;  "AC" = 100 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

$excel.Calculate()
